$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5757.4644
$ws.Range("I106").Value = 5859.593
$ws.Range("K106").Value = 5859.593
$ws.Range("M106").Value = -5228.593
$ws.Range("H110").Value = 32702
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H113").Value = 3420
$ws.Range("I113").Value = 3493.3333
$ws.Range("J113").Value = 3200
$ws.Range("K113").Value = 3493.3333
$ws.Range("L113").Value = 3200
$ws.Range("M113").Value = -239.3332999999998
$ws.Range("N113").Value = -9708
$ws.Range("H132").Value = 8556036
$ws.Range("I132").Value = 11914045
$ws.Range("J132").Value = 8374.817999999999
$ws.Range("K132").Value = 35742135
$ws.Range("L132").Value = 25124.454
$ws.Range("M132").Value = -35739605
$ws.Range("N132").Value = -30184.454
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H138").Value = 558988
$ws.Range("I138").Value = 1204.2941
$ws.Range("J138").Value = 714435.9
$ws.Range("K138").Value = 3612.8823
$ws.Range("L138").Value = 2143307.7
$ws.Range("M138").Value = 1527.1177
$ws.Range("N138").Value = -2153587.7
$ws.Range("H140").Value = 33028.57
$ws.Range("J140").Value = 33028.57
$ws.Range("L140").Value = 33028.57
$ws.Range("N140").Value = -43388.57

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4787.7163
$ws.Range("I32").Value = 4668.4917
$ws.Range("K32").Value = 4668.4917
$ws.Range("M32").Value = -4381.4917
$ws.Range("H33").Value = 100000000
$ws.Range("I33").Value = 100000000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 100000000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -99999671
$ws.Range("N33").ClearContents()
$ws.Range("H35").Value = 901
$ws.Range("I35").Value = 901
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 901
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -495
$ws.Range("N35").ClearContents()
$ws.Range("H36").Value = 1000
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H39").Value = 3900
$ws.Range("I39").Value = 3900
$ws.Range("K39").Value = 3900
$ws.Range("M39").Value = -3380
$ws.Range("H41").Value = 7351.2
$ws.Range("I41").Value = 5689
$ws.Range("J41").Value = 14000
$ws.Range("K41").Value = 5689
$ws.Range("L41").Value = 14000
$ws.Range("M41").Value = -5275
$ws.Range("N41").Value = -14828
$ws.Range("H74").Value = 1716.9565
$ws.Range("I74").Value = 979.4
$ws.Range("K74").Value = 979.4
$ws.Range("M74").Value = -105.4
$ws.Range("H77").Value = 1716.9565
$ws.Range("I77").Value = 979.4
$ws.Range("K77").Value = 4897
$ws.Range("M77").Value = -529
$ws.Range("H122").Value = 2759.2
$ws.Range("I122").Value = 2784.8
$ws.Range("J122").Value = 2733.6
$ws.Range("K122").Value = 8354.400000000001
$ws.Range("L122").Value = 8200.799999999999
$ws.Range("M122").Value = -5904.400000000001
$ws.Range("N122").Value = -13100.8
$ws.Range("H139").Value = 28593.334
$ws.Range("J139").Value = 28593.334
$ws.Range("L139").Value = 28593.334
$ws.Range("N139").Value = -38873.334
$ws.Range("H140").Value = 33357.145
$ws.Range("J140").Value = 33357.145
$ws.Range("L140").Value = 33357.145
$ws.Range("N140").Value = -43717.145
$ws.Range("H141").Value = 32280.777
$ws.Range("J141").Value = 32280.777
$ws.Range("L141").Value = 32280.777
$ws.Range("N141").Value = -42640.777

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1471.9333
$ws.Range("I134").Value = 1189.2727
$ws.Range("K134").Value = 3567.8181
$ws.Range("M134").Value = -1032.8181
$ws.Range("H140").Value = 23917.5
$ws.Range("J140").Value = 23917.5
$ws.Range("L140").Value = 23917.5
$ws.Range("N140").Value = -34277.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1356.6459
$ws.Range("I31").Value = 1211.8182
$ws.Range("J31").Value = 2949.75
$ws.Range("K31").Value = 1211.8182
$ws.Range("L31").Value = 2949.75
$ws.Range("M31").Value = -916.8181999999999
$ws.Range("N31").Value = -3539.75
$ws.Range("H34").Value = 1356.6459
$ws.Range("I34").Value = 1211.8182
$ws.Range("J34").Value = 2949.75
$ws.Range("K34").Value = 1211.8182
$ws.Range("L34").Value = 2949.75
$ws.Range("M34").Value = -1009.8182
$ws.Range("N34").Value = -3353.75
$ws.Range("H62").Value = 11113360
$ws.Range("I62").Value = 2328.5715
$ws.Range("J62").Value = 50001970
$ws.Range("K62").Value = 2328.5715
$ws.Range("L62").Value = 50001970
$ws.Range("M62").Value = -1704.5715
$ws.Range("N62").Value = -50003218
$ws.Range("H65").Value = 11113360
$ws.Range("I65").Value = 2328.5715
$ws.Range("J65").Value = 50001970
$ws.Range("K65").Value = 11642.8575
$ws.Range("L65").Value = 250009850
$ws.Range("M65").Value = -8522.8575
$ws.Range("N65").Value = -250016090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 863.8461
$ws.Range("I2").Value = 63.22222
$ws.Range("K2").Value = 379.33332
$ws.Range("M2").Value = -266.33332
$ws.Range("H131").Value = 15152402
$ws.Range("J131").Value = 1015.28845
$ws.Range("L131").Value = 3045.86535
$ws.Range("N131").Value = -13125.86535
$ws.Range("H139").Value = 1951.4878
$ws.Range("I139").Value = 2093
$ws.Range("J139").Value = 1678.5714
$ws.Range("K139").Value = 6279
$ws.Range("L139").Value = 5035.7142
$ws.Range("M139").Value = -1139
$ws.Range("N139").Value = -15315.7142
$ws.Range("H140").Value = 24042.086
$ws.Range("I140").Value = 49476.43
$ws.Range("J140").Value = 3498.9614
$ws.Range("K140").Value = 148429.29
$ws.Range("L140").Value = 10496.8842
$ws.Range("M140").Value = -143249.29
$ws.Range("N140").Value = -20856.8842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1231.1333
$ws.Range("I102").Value = 1544.8
$ws.Range("J102").Value = 1074.3
$ws.Range("K102").Value = 1544.8
$ws.Range("L102").Value = 1074.3
$ws.Range("M102").Value = 77.20000000000005
$ws.Range("N102").Value = -4318.3
$ws.Range("H122").Value = 2798.25
$ws.Range("I122").Value = 3345.5715
$ws.Range("J122").Value = 2032
$ws.Range("K122").Value = 10036.7145
$ws.Range("L122").Value = 6096
$ws.Range("M122").Value = -7586.7145
$ws.Range("N122").Value = -10996
$ws.Range("H136").Value = 9395.625
$ws.Range("J136").Value = 9395.625
$ws.Range("L136").Value = 28186.875
$ws.Range("N136").Value = -33286.875
$ws.Range("H141").Value = 27242.857
$ws.Range("J141").Value = 27242.857
$ws.Range("L141").Value = 27242.857
$ws.Range("N141").Value = -37602.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2213.4285
$ws.Range("I7").Value = 2213.4285
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2213.4285
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2101.4285
$ws.Range("N7").ClearContents()
$ws.Range("H40").Value = 3569
$ws.Range("I40").Value = 2192.75
$ws.Range("K40").Value = 2192.75
$ws.Range("M40").Value = -2056.75
$ws.Range("H68").Value = 2028.2354
$ws.Range("I68").Value = 1998.75
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 1998.75
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -1249.75
$ws.Range("N68").Value = -3998
$ws.Range("H71").Value = 2028.2354
$ws.Range("I71").Value = 1998.75
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 9993.75
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -6249.75
$ws.Range("N71").Value = -19988
$ws.Range("H122").Value = 35716170
$ws.Range("I122").Value = 50001836
$ws.Range("J122").Value = 2002.5
$ws.Range("K122").Value = 150005508
$ws.Range("L122").Value = 6007.5
$ws.Range("M122").Value = -150003058
$ws.Range("N122").Value = -10907.5
$ws.Range("H126").Value = 2213.4285
$ws.Range("I126").Value = 2213.4285
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6640.2855
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4170.2855
$ws.Range("N126").ClearContents()
$ws.Range("H138").Value = 34600
$ws.Range("J138").Value = 34600
$ws.Range("L138").Value = 34600
$ws.Range("N138").Value = -44880
$ws.Range("H140").Value = 51320.9
$ws.Range("J140").Value = 51320.9
$ws.Range("L140").Value = 51320.9
$ws.Range("N140").Value = -61680.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 125004616
$ws.Range("I126").Value = 142861000
$ws.Range("K126").Value = 428583000
$ws.Range("M126").Value = -428580530
$ws.Range("H140").Value = 28179.273
$ws.Range("J140").Value = 29997.2
$ws.Range("L140").Value = 29997.2
$ws.Range("N140").Value = -40357.2
